$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$tbl = $s.Shapes.Item(3).Table
$tbl.ApplyStyle("{9C86746E-C20C-45AA-8946-1B3ABE58F2AF}")
